# Applies the BBPPRTY "Max Potential Capacity (MW)" -> "% of capacity
# available for retrofit" header rename, and zeroes out the 2028-2050
# (I19:AE19) "hard coal w CCS" retrofit-allowed flags (previously all 1s)
# while giving them the same integer number format used by the other
# data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BBPPRTY")

# Rename the header cell that used to read "Max Potential Capacity (MW)".
$ws.Range("A1").Value = "% of capacity available for retrofit"

# 2028 (column I) through 2050 (column AE) on the "hard coal w CCS" row (19)
# go from allowed (1) to not allowed (0), formatted like the rest of the
# boolean-flag rows (integer "0" number format).
$rng = $ws.Range("I19:AE19")
$rng.Value = 0
$rng.NumberFormat = "0"

# Move the lingering selection on the BBPPRTY sheet from B29 to A2, then
# restore "About" as the active/visible tab (it was tabSelected before the
# edit and should remain so).
$ws.Range("A2").Select()
$about = $wb.Worksheets.Item("About")
$about.Activate()
